{"js": "// Split the paragraph describing the new pizzeria into two paragraphs.\n//\n// Before (single paragraph, single run of text):\n//   \"...offrendo una vasta scelta di pizze, incluse proposte gourmet realizzate\n//    con prodotti forniti da aziende agricole locali. Questo modello non solo\n//    instaura un legame con il territorio, ma promuove anche un approccio\n//    sostenibile e a chilometro zero, utilizzando materie prime biologiche e\n//    di alta qualit\u00e0.\"\n//\n// After (two paragraphs): the sentence starting at \"Questo modello\" moves\n// into its own, new paragraph; the first paragraph keeps the trailing space\n// after \"locali. \".\n\nconst body = context.document.body;\n\nconst tailText =\n  \"Questo modello non solo instaura un legame con il territorio, ma \" +\n  \"promuove anche un approccio sostenibile e a chilometro zero, utilizzando \" +\n  \"materie prime biologiche e di alta qualit\u00e0.\";\n\nconst searchResults = body.search(tailText, { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Could not find the sentence that starts the new paragraph.\");\n}\n\nconst tailRange = searchResults.items[0];\n\n// `Range.insertParagraph` operates at the level of the whole paragraph that\n// contains the range, so first append the tail sentence as a brand-new\n// paragraph right after the current one...\ntailRange.insertParagraph(tailText, Word.InsertLocation.after);\nawait context.sync();\n\n// ...then strip the original (now duplicated) copy of that sentence out of\n// the first paragraph, leaving only the text up to and including \"locali. \".\nconst dupResults = body.search(tailText, { matchCase: true });\ndupResults.load(\"items\");\nawait context.sync();\n\ndupResults.items[0].insertText(\"\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Split the paragraph describing the new pizzeria into two paragraphs.\n#\n# Before (single paragraph, single run of text):\n#   \"...offrendo una vasta scelta di pizze, incluse proposte gourmet realizzate\n#    con prodotti forniti da aziende agricole locali. Questo modello non solo\n#    instaura un legame con il territorio, ma promuove anche un approccio\n#    sostenibile e a chilometro zero, utilizzando materie prime biologiche e\n#    di alta qualit\u00e0.\"\n#\n# After (two paragraphs): the sentence starting at \"Questo modello\" moves\n# into its own, new paragraph; the first paragraph keeps the trailing space\n# after \"locali. \".\n\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = \"Questo modello non solo instaura\"\n$find.MatchCase = $true\n$found = $find.Execute()\n\nif ($found) {\n    # Insert a paragraph break immediately before the matched text, which\n    # pushes it (and everything after it) into a new paragraph.\n    $rng.InsertParagraphBefore()\n}\n"}
